$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1845.4166
$ws.Range("I112").Value = 300
$ws.Range("J112").Value = 1985.909
$ws.Range("K112").Value = 900
$ws.Range("L112").Value = 5957.727000000001
$ws.Range("M112").Value = 208
$ws.Range("N112").Value = -8173.727000000001

$ws.Range("H116").Value = 2676.4
$ws.Range("I116").Value = 2422.7273
$ws.Range("J116").Value = 3374
$ws.Range("K116").Value = 2422.7273
$ws.Range("L116").Value = 3374
$ws.Range("M116").Value = 1019.2727
$ws.Range("N116").Value = -10258

$ws.Range("H129").Value = 849.73914
$ws.Range("J129").Value = 943.2820400000001
$ws.Range("L129").Value = 2829.84612
$ws.Range("N129").Value = -12829.84612

$ws.Range("H137").Value = 3208.2368
$ws.Range("I137").Value = 1767.0416
$ws.Range("J137").Value = 5678.857
$ws.Range("K137").Value = 5301.1248
$ws.Range("L137").Value = 17036.571
$ws.Range("M137").Value = -2751.1248
$ws.Range("N137").Value = -22136.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6859.6924
$ws.Range("I32").Value = 5543.8613
$ws.Range("K32").Value = 5543.8613
$ws.Range("M32").Value = -5256.8613

$ws.Range("H61").Value = 10783.833
$ws.Range("J61").Value = 15471
$ws.Range("L61").Value = 15471
$ws.Range("N61").Value = -15895

$ws.Range("H74").Value = 6431.522
$ws.Range("I74").Value = 8400.6875
$ws.Range("J74").Value = 1930.5714
$ws.Range("K74").Value = 8400.6875
$ws.Range("L74").Value = 1930.5714
$ws.Range("M74").Value = -7526.6875
$ws.Range("N74").Value = -3678.5714

$ws.Range("H77").Value = 6431.522
$ws.Range("I77").Value = 8400.6875
$ws.Range("J77").Value = 1930.5714
$ws.Range("K77").Value = 42003.4375
$ws.Range("L77").Value = 9652.857
$ws.Range("M77").Value = -37635.4375
$ws.Range("N77").Value = -18388.857

$ws.Range("H132").Value = 3983.1853
$ws.Range("I132").Value = 2376.5
$ws.Range("K132").Value = 7129.5
$ws.Range("M132").Value = -4599.5

$ws.Range("H136").Value = 10783.833
$ws.Range("J136").Value = 15471
$ws.Range("L136").Value = 46413
$ws.Range("N136").Value = -51513

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2909
$ws.Range("I134").Value = 2319.8572
$ws.Range("K134").Value = 6959.571599999999
$ws.Range("M134").Value = -4424.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3137999
$ws.Range("I58").Value = 5052622
$ws.Range("K58").Value = 5052622
$ws.Range("M58").Value = -5052419

$ws.Range("H132").Value = 2899.4583
$ws.Range("I132").Value = 2286.8125
$ws.Range("J132").Value = 4124.75
$ws.Range("K132").Value = 6860.4375
$ws.Range("L132").Value = 12374.25
$ws.Range("M132").Value = -4330.4375
$ws.Range("N132").Value = -17434.25

$ws.Range("H134").Value = 3637.3547
$ws.Range("I134").Value = 2798
$ws.Range("J134").Value = 4799.5386
$ws.Range("K134").Value = 8394
$ws.Range("L134").Value = 14398.6158
$ws.Range("M134").Value = -5859
$ws.Range("N134").Value = -19468.6158

$ws.Range("H136").Value = 3137999
$ws.Range("I136").Value = 5052622
$ws.Range("K136").Value = 15157866
$ws.Range("M136").Value = -15155316

$ws.Range("H141").Value = 42434.125
$ws.Range("I141").Value = 37648
$ws.Range("J141").Value = 44029.5
$ws.Range("K141").Value = 37648
$ws.Range("L141").Value = 44029.5
$ws.Range("M141").Value = -32468
$ws.Range("N141").Value = -54389.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4765544
$ws.Range("I5").Value = 454
$ws.Range("J5").Value = 11118997
$ws.Range("K5").Value = 1362
$ws.Range("L5").Value = 33356991
$ws.Range("M5").Value = -1250
$ws.Range("N5").Value = -33357215

$ws.Range("H46").Value = 2497.5293
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 2634.875
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 7904.625
$ws.Range("M46").Value = -809
$ws.Range("N46").Value = -8086.625

$ws.Range("H113").Value = 677.13
$ws.Range("J113").Value = 640
$ws.Range("L113").Value = 1920
$ws.Range("N113").Value = -6260

$ws.Range("H135").Value = 4765544
$ws.Range("I135").Value = 454
$ws.Range("J135").Value = 11118997
$ws.Range("K135").Value = 4086
$ws.Range("L135").Value = 100070973
$ws.Range("M135").Value = -1551
$ws.Range("N135").Value = -100076043

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3343.6667
$ws.Range("I132").Value = 2925.3333
$ws.Range("J132").Value = 3762
$ws.Range("K132").Value = 8775.999899999999
$ws.Range("L132").Value = 11286
$ws.Range("M132").Value = -6245.999899999999
$ws.Range("N132").Value = -16346

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H82").Value = 1883.2858
$ws.Range("I82").Value = 1149
$ws.Range("J82").Value = 2005.6666
$ws.Range("K82").Value = 1149
$ws.Range("L82").Value = 2005.6666
$ws.Range("M82").Value = -788
$ws.Range("N82").Value = -2727.6666

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H85").Value = 1883.2858
$ws.Range("I85").Value = 1149
$ws.Range("J85").Value = 2005.6666
$ws.Range("K85").Value = 1149
$ws.Range("L85").Value = 2005.6666
$ws.Range("M85").Value = 99
$ws.Range("N85").Value = -4501.6666

$ws.Range("H136").Value = 6436.579
$ws.Range("I136").Value = 5136.579
$ws.Range("K136").Value = 15409.737
$ws.Range("M136").Value = -12859.737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1829.2307
$ws.Range("I132").Value = 880.5714
$ws.Range("J132").Value = 2936
$ws.Range("K132").Value = 2641.7142
$ws.Range("L132").Value = 8808
$ws.Range("M132").Value = -111.7142000000003
$ws.Range("N132").Value = -13868

$ws.Range("H136").Value = 4040.2354
$ws.Range("I136").Value = 3506.611
$ws.Range("J136").Value = 4640.5625
$ws.Range("K136").Value = 10519.833
$ws.Range("L136").Value = 13921.6875
$ws.Range("M136").Value = -7969.832999999999
$ws.Range("N136").Value = -19021.6875

